$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text number format on the numeric-looking columns that must stay as
# literal text (Price, Volume%, Hora), matching the original inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "307.42"
$ws.Range("E2").Value = "0.99%"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "36.43"
$ws.Range("E3").Value = "1.67%"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "5.051"
$ws.Range("E4").Value = "-0.38%"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.08103"
$ws.Range("E5").Value = "0.62%"
$ws.Range("G5").Value = "18"
$ws.Range("D6").Value = "2.126"
$ws.Range("E6").Value = "10.18%"
$ws.Range("G6").Value = "18"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "7.842"
$ws.Range("E7").Value = "-0.02%"
$ws.Range("G7").Value = "18"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9306"
$ws.Range("E8").Value = "0.03%"
$ws.Range("G8").Value = "18"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1454"
$ws.Range("E9").Value = "14.60%"
$ws.Range("G9").Value = "18"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1923"
$ws.Range("E10").Value = "-0.06%"
$ws.Range("G10").Value = "18"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09116"
$ws.Range("E11").Value = "-0.51%"
$ws.Range("G11").Value = "18"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03453"
$ws.Range("E12").Value = "-0.74%"
$ws.Range("G12").Value = "18"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09916"
$ws.Range("E13").Value = "0.18%"
$ws.Range("G13").Value = "18"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001409"
$ws.Range("E14").Value = "-0.70%"
$ws.Range("G14").Value = "18"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006302"
$ws.Range("E15").Value = "-5.34%"
$ws.Range("G15").Value = "18"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.840"
$ws.Range("E16").Value = "6.29%"
$ws.Range("G16").Value = "18"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.145"
$ws.Range("E17").Value = "-0.10%"
$ws.Range("G17").Value = "18"
$ws.Range("D18").Value = "3.485"
$ws.Range("E18").Value = "10.01%"
$ws.Range("G18").Value = "18"
$ws.Range("D19").Value = "0.3458"
$ws.Range("E19").Value = "1.07%"
$ws.Range("G19").Value = "18"
$ws.Range("D20").Value = "0.1282"
$ws.Range("E20").Value = "-4.06%"
$ws.Range("G20").Value = "18"
$ws.Range("D21").Value = "4.799"
$ws.Range("E21").Value = "-7.28%"
$ws.Range("G21").Value = "18"
$ws.Range("D22").Value = "0.2335"
$ws.Range("E22").Value = "-7.79%"
$ws.Range("G22").Value = "18"
$ws.Range("D23").Value = "0.04356"
$ws.Range("E23").Value = "-1.21%"
$ws.Range("G23").Value = "18"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").Value = "-0.76%"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.004919"
$ws.Range("E25").Value = "4.16%"
$ws.Range("G25").Value = "18"
$ws.Range("G26").Value = "18"
$ws.Range("D27").Value = "0.0001297"
$ws.Range("E27").Value = "-0.50%"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("G38").Value = "18"
$ws.Range("D39").Value = "0.02027"
$ws.Range("E39").Value = "1.50%"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.05206"
$ws.Range("E40").Value = "0.97%"
$ws.Range("G40").Value = "18"
$ws.Range("D41").Value = "0.007533"
$ws.Range("E41").Value = "-1.08%"
$ws.Range("G41").Value = "18"
$ws.Range("D42").Value = "0.01009"
$ws.Range("E42").Value = "0.34%"
$ws.Range("G42").Value = "18"
$ws.Range("D43").Value = "0.1371"
$ws.Range("E43").Value = "0.34%"
$ws.Range("G43").Value = "18"
$ws.Range("D44").Value = "0.002144"
$ws.Range("E44").Value = "1.86%"
$ws.Range("G44").Value = "18"
$ws.Range("D45").Value = "0.009962"
$ws.Range("E45").Value = "-6.81%"
$ws.Range("G45").Value = "18"
$ws.Range("D46").Value = "0.00006273"
$ws.Range("E46").Value = "-1.11%"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.35%"
$ws.Range("G47").Value = "18"
$ws.Range("D48").Value = "64.90"
$ws.Range("E48").Value = "-0.50%"
$ws.Range("G48").Value = "18"
$ws.Range("D49").Value = "0.001247"
$ws.Range("E49").Value = "-22.15%"
$ws.Range("G49").Value = "18"
$ws.Range("D50").Value = "0.00002094"
$ws.Range("E50").Value = "-0.35%"
$ws.Range("G50").Value = "18"
$ws.Range("D51").Value = "0.0001995"
$ws.Range("E51").Value = "-0.35%"
$ws.Range("G51").Value = "18"
